$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CC / 1043312708 / JUAN LUIS CASTELLANO MARTINEZ / 1902
$ws.Range("C16").Value = "1043312708"
$ws.Range("D16").Value = "JUAN LUIS CASTELLANO MARTINEZ"
$ws.Range("E16").Value = "1902"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

# Row 17: CC / 1043312708 / JUAN LUIS CASTELLANO MARTINEZ / 1903
$ws.Range("C17").Value = "1043312708"
$ws.Range("D17").Value = "JUAN LUIS CASTELLANO MARTINEZ"
$ws.Range("E17").Value = "1903"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 828116

# Row 18: CC / 1043312708 / JUAN LUIS CASTELLANO MARTINEZ / 2007
$ws.Range("C18").Value = "1043312708"
$ws.Range("D18").Value = "JUAN LUIS CASTELLANO MARTINEZ"
$ws.Range("E18").Value = "2007"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 877803

# Row 19: CC / 1047370395 / OCTAVIO ENRIQUE RAUDALES GARRIDO / 2007
$ws.Range("C19").Value = "1047370395"
$ws.Range("D19").Value = "OCTAVIO ENRIQUE RAUDALES GARRIDO"
$ws.Range("E19").Value = "2007"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = 877803

# Row 20: CC / 7921854 / ROBERTO CARLOS ATENCIO CUELLO / 2008
$ws.Range("C20").Value = "7921854"
$ws.Range("D20").Value = "ROBERTO CARLOS ATENCIO CUELLO"
$ws.Range("E20").Value = "2008"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 877803

# Row 21: CC / 1047370395 / OCTAVIO ENRIQUE RAUDALES GARRIDO / 2008
$ws.Range("C21").Value = "1047370395"
$ws.Range("D21").Value = "OCTAVIO ENRIQUE RAUDALES GARRIDO"
$ws.Range("E21").Value = "2008"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 877803
